$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Kyshawn George -> De'Aaron Fox (PG, Sacramento Kings)
$ws.Range("A4").Value = "De'Aaron Fox"
$ws.Range("B4").Value = "PG"
$ws.Range("C4").Value = "Sacramento Kings"

# Row 5: De'Aaron Fox -> Josh Giddey (PG,SG,SF, Chicago Bulls)
$ws.Range("A5").Value = "Josh Giddey"
$ws.Range("B5").Value = "PG,SG,SF"
$ws.Range("C5").Value = "Chicago Bulls"

# Row 12: Josh Giddey -> Nikola Vucevic (PF,C, Chicago Bulls)
$ws.Range("A12").Value = "Nikola Vucevic"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Chicago Bulls"

# Row 13: Nikola Vucevic -> Isaiah Stewart (PF,C stays, Detroit Pistons)
$ws.Range("A13").Value = "Isaiah Stewart"
$ws.Range("C13").Value = "Detroit Pistons"
